$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 292. Excel will shift all existing
# rows 292..351 down to 293..352, which matches the target diff (every
# previously-existing row keeps its data but moves down by one row).
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new weekly record.
$ws.Cells.Item(292, 1).Value = 7
$ws.Cells.Item(292, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(292, 3).Value = "Ñuble"
$ws.Cells.Item(292, 4).Value = (Get-Date -Year 2023 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(292, 5).Value = 16
$ws.Cells.Item(292, 6).Value = 100112045
$ws.Cells.Item(292, 7).Value = "Zapallo"
$ws.Cells.Item(292, 8).Value = "Camote"
$ws.Cells.Item(292, 9).Value = "1a (guarda)"
$ws.Cells.Item(292, 10).Value = 250
$ws.Cells.Item(292, 11).Value = 650
$ws.Cells.Item(292, 12).Value = 650
$ws.Cells.Item(292, 13).Value = 650
$ws.Cells.Item(292, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(292, 15).Value = "Región del Maule"
$ws.Cells.Item(292, 16).Value = 650
$ws.Cells.Item(292, 17).Value = 1
$ws.Cells.Item(292, 18).Value = "Hortaliza"

# Give the new date cell the same date number-format used by the rest of
# column D.
$ws.Cells.Item(292, 4).NumberFormat = $ws.Cells.Item(293, 4).NumberFormat
